# Updated latest Guinea master data.
# Rewrites the doc_category sheet to the new MOSIP export layout: adds the
# standard audit columns (cr_by, cr_dtimes, upd_by, upd_dtimes, is_deleted,
# del_dtimes), re-types is_active/is_deleted as real booleans, fixes a couple
# of mis-encoded accented French strings, and reorders the POC/POE rows to
# match the latest export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row -----------------------------------------------------
$ws.Range("A1").Value = "code"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "descr"
$ws.Range("D1").Value = "lang_code"
$ws.Range("E1").Value = "is_active"
$ws.Range("F1").Value = "cr_by"
$ws.Range("G1").Value = "cr_dtimes"
$ws.Range("H1").Value = "upd_by"
$ws.Range("I1").Value = "upd_dtimes"
$ws.Range("J1").Value = "is_deleted"
$ws.Range("K1").Value = "del_dtimes"

# ---- Shared constant values ------------------------------------------
$crBy    = "superadmin"
$crDt    = 45079.576954143522
$nullTxt = "NULL"
$updBy   = "chithara27"

# ---- Data rows --------------------------------------------------------
# Row 2: POA
$ws.Range("A2").Value = "POA"
$ws.Range("B2").Value = "Justificatif de domicile"
$ws.Range("C2").Value = "Justificatif de domicile"
$ws.Range("D2").Value = "fra"
$ws.Range("E2").Value = $true
$ws.Range("F2").Value = $crBy
$ws.Range("G2").Value = $crDt
$ws.Range("G2").NumberFormat = "mm:ss.0"
$ws.Range("H2").Value = $nullTxt
$ws.Range("I2").Value = $nullTxt
$ws.Range("J2").Value = $false
$ws.Range("K2").Value = $nullTxt

# Row 3: POI
$ws.Range("A3").Value = "POI"
$ws.Range("B3").Value = "Justificatif d'identitÃ©"
$ws.Range("C3").Value = "Justificatif d'identitÃ©"
$ws.Range("D3").Value = "fra"
$ws.Range("E3").Value = $true
$ws.Range("F3").Value = $crBy
$ws.Range("G3").Value = $crDt
$ws.Range("G3").NumberFormat = "mm:ss.0"
$ws.Range("H3").Value = $nullTxt
$ws.Range("I3").Value = $nullTxt
$ws.Range("J3").Value = $false
$ws.Range("K3").Value = $nullTxt

# Row 4: POR
$ws.Range("A4").Value = "POR"
$ws.Range("B4").Value = "Justificatif de lien de parentÃ©"
$ws.Range("C4").Value = "Justificatif de lien de parentÃ©"
$ws.Range("D4").Value = "fra"
$ws.Range("E4").Value = $true
$ws.Range("F4").Value = $crBy
$ws.Range("G4").Value = $crDt
$ws.Range("G4").NumberFormat = "mm:ss.0"
$ws.Range("H4").Value = $nullTxt
$ws.Range("I4").Value = $nullTxt
$ws.Range("J4").Value = $false
$ws.Range("K4").Value = $nullTxt

# Row 5: POB
$ws.Range("A5").Value = "POB"
$ws.Range("B5").Value = "Justificatif de date de naissance"
$ws.Range("C5").Value = "Justificatif de date de naissance"
$ws.Range("D5").Value = "fra"
$ws.Range("E5").Value = $true
$ws.Range("F5").Value = $crBy
$ws.Range("G5").Value = $crDt
$ws.Range("G5").NumberFormat = "mm:ss.0"
$ws.Range("H5").Value = $nullTxt
$ws.Range("I5").Value = $nullTxt
$ws.Range("J5").Value = $false
$ws.Range("K5").Value = $nullTxt

# Row 6: POC (moved up from old row 7, now deleted/deactivated)
$ws.Range("A6").Value = "POC"
$ws.Range("B6").Value = "Justificatif de consentement"
$ws.Range("C6").Value = "Justificatif de consentement"
$ws.Range("D6").Value = "fra"
$ws.Range("E6").Value = $false
$ws.Range("F6").Value = $crBy
$ws.Range("G6").Value = $crDt
$ws.Range("G6").NumberFormat = "mm:ss.0"
$ws.Range("H6").Value = $updBy
$ws.Range("I6").Value = 45091.286432222223
$ws.Range("I6").NumberFormat = "mm:ss.0"
$ws.Range("J6").Value = $false
$ws.Range("K6").Value = $nullTxt

# Row 7: POE (moved down from old row 6, now deleted/deactivated)
$ws.Range("A7").Value = "POE"
$ws.Range("B7").Value = "Justificatif d'exception biomÃ©trique"
$ws.Range("C7").Value = "Justificatif d'exception biomÃ©trique"
$ws.Range("D7").Value = "fra"
$ws.Range("E7").Value = $false
$ws.Range("F7").Value = $crBy
$ws.Range("G7").Value = $crDt
$ws.Range("G7").NumberFormat = "mm:ss.0"
$ws.Range("H7").Value = $updBy
$ws.Range("I7").Value = 45091.286544212962
$ws.Range("I7").NumberFormat = "mm:ss.0"
$ws.Range("J7").Value = $false
$ws.Range("K7").Value = $nullTxt

# ---- Selection, matching the refreshed export ------------------------
$ws.Range("E14").Select()
